$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '39.705.90'
$ws.Range("E2").Value = '  +2.34%  '
$ws.Range("D3").Value = '2.164.65'
$ws.Range("E3").Value = '  +2.83%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '228.03'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  +0.11%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.629'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = '  +2.03%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '63.56'
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = '  +1.51%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("E9").Value = '  +0.91%  '
$ws.Range("E10").Value = '  +0.62%  '
$ws.Range("E11").Value = '  +0.12%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '16.06'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = '  +1.56%  '
$ws.Range("D13").Value = '2.484.49'
$ws.Range("E13").Value = '  +2.82%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '22.04'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = '  -0.02%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.810'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = '  +0.00%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.50'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = '  -0.71%  '
$ws.Range("D17").Value = '2.163.75'
$ws.Range("E17").Value = '  +2.91%  '
$ws.Range("D18").Value = '39.615.96'
$ws.Range("E18").Value = '  +2.10%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '71.91'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = '  +0.52%  '
$ws.Range("E20").Value = '  +0.21%  '
$ws.Range("E21").Value = '  -0.05%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '228.14'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = '  -0.05%  '
$ws.Range("E23").Value = '  +0.07%  '
$ws.Range("E24").Value = '  +3.39%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.32'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = '  -1.92%  '
$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '172.52'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = '  +0.09%  '
$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.49'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = '  -1.78%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.139'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = '  +0.75%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.77'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = '  +2.19%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.41'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = '  +0.16%  '
$ws.Range("E31").Value = '  +4.34%  '
$ws.Range("E32").Value = '  +1.61%  '
$ws.Range("E33").Value = '  +0.67%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.69'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = '  -1.43%  '
$ws.Range("E35").Value = '  -3.04%  '
$ws.Range("E36").Value = '  +0.10%  '
$ws.Range("E37").Value = '  +0.66%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.64'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = '  +3.21%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.999'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = '  +0.06%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.82'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = '  +15.17%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '102.13'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = '  -0.30%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0227'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = '  -0.24%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '17.63'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = '  -2.51%  '
$ws.Range("B44").Value = 'TrustWalletToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.22'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = '  +1.24%  '
$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").Value = '1.511.41'
$ws.Range("E45").Value = '  -1.00%  '
$ws.Range("E46").Value = '  +0.80%  '
$ws.Range("B47").Value = 'ARBITRUM'
$ws.Range("C47").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.10'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = '  +1.97%  '
$ws.Range("B48").Value = 'HuobiToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.80'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = '  -0.21%  '
$ws.Range("E49").Value = '  -0.43%  '
$ws.Range("E50").Value = '  +1.21%  '
$ws.Range("D51").Value = '2.368.54'
